$d = $word.ActiveDocument

# --- Change 1: merge "190к" + " слов" into a single run "190к слов" ---
$d.Content.Find.Execute("190к слов", $true, $false, $false, $false, $false, $true, 1, $false, "190к слов", 2) | Out-Null

# --- Change 2: merge the run sequence in the "Разговорник" paragraph ---
$old2 = " содержит все необходимые слова и фразы, которые могут Вам понадобиться во время вашей работы поездки или путешествия. База разговорника содержит "
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

# --- Change 3: merge the tag-list run sequence around "таджикско-таджикский (толковый словарь)" ---
$old3 = ", переводчик, перевод, словарь, таджикский, русский, таджикско-русский, русско-таджикский, таджикско-таджикский (толковый словарь), "
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2) | Out-Null

# --- Change 4: remove the LISTNUM field code (fldChar begin / instrText / fldChar end runs) ---
for ($i = $d.Fields.Count; $i -ge 1; $i--) {
    $d.Fields.Item($i).Delete()
}
